# Update column G ("K" - strikeouts) values for rows 2-41 on the active sheet.
# These values were regenerated to use K instead of Strike# (old Strike count data
# was replaced by true strikeout totals per game).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 2
    4  = 5
    5  = 5
    6  = 10
    7  = 3
    8  = 5
    9  = 5
    10 = 7
    11 = 8
    12 = 9
    13 = 4
    14 = 5
    15 = 6
    16 = 7
    17 = 5
    18 = 5
    19 = 3
    20 = 13
    21 = 3
    22 = 12
    23 = 6
    24 = 5
    25 = 5
    26 = 5
    27 = 5
    28 = 10
    29 = 8
    30 = 6
    31 = 4
    32 = 10
    33 = 8
    34 = 11
    35 = 7
    36 = 3
    37 = 6
    38 = 3
    39 = 3
    40 = 5
    41 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
